# Update row 8 (year 2025) metrics in the active worksheet of the workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1214
$ws.Range("D8").Value = 198
$ws.Range("E8").Value = 1016
$ws.Range("F8").Value = 8.1214109926169
$ws.Range("G8").Value = 83.69028006589787
$ws.Range("H8").Value = 16.30971993410214
